$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 878.2619
$ws.Range("J17").Value = 878.2619
$ws.Range("L17").Value = 2634.7857
$ws.Range("N17").Value = -2970.7857
$ws.Range("H51").Value = 9949.75
$ws.Range("I51").Value = 10000
$ws.Range("J51").Value = 9944.166999999999
$ws.Range("K51").Value = 10000
$ws.Range("L51").Value = 9944.166999999999
$ws.Range("M51").Value = -9516
$ws.Range("N51").Value = -10912.167
$ws.Range("H62").Value = 5000
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 5000
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 5000
$ws.Range("M62").Value = ""
$ws.Range("N62").Value = -6248
$ws.Range("H65").Value = 5000
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 5000
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 25000
$ws.Range("M65").Value = ""
$ws.Range("N65").Value = -31240
$ws.Range("H86").Value = 2300.5715
$ws.Range("I86").Value = 2329.4546
$ws.Range("J86").Value = 2194.6667
$ws.Range("K86").Value = 2329.4546
$ws.Range("L86").Value = 2194.6667
$ws.Range("M86").Value = -1206.4546
$ws.Range("N86").Value = -4440.6667
$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").Value = ""
$ws.Range("H89").Value = 2300.5715
$ws.Range("I89").Value = 2329.4546
$ws.Range("J89").Value = 2194.6667
$ws.Range("K89").Value = 11647.273
$ws.Range("L89").Value = 10973.3335
$ws.Range("M89").Value = -6031.273000000001
$ws.Range("N89").Value = -22205.3335
$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").Value = ""
$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").Value = ""
$ws.Range("H111").Value = 566.3333
$ws.Range("I111").Value = 349.5
$ws.Range("K111").Value = 1048.5
$ws.Range("M111").Value = 2018.5
$ws.Range("H138").Value = 23257756
$ws.Range("I138").Value = 1877.6471
$ws.Range("J138").Value = 38463520
$ws.Range("K138").Value = 5632.9413
$ws.Range("L138").Value = 115390560
$ws.Range("M138").Value = -492.9412999999995
$ws.Range("N138").Value = -115400840
$ws.Range("H141").Value = 2704.12
$ws.Range("I141").Value = 2386.2856
$ws.Range("J141").Value = 4372.75
$ws.Range("K141").Value = 7158.8568
$ws.Range("L141").Value = 13118.25
$ws.Range("M141").Value = -1978.8568
$ws.Range("N141").Value = -23478.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 3764.0476
$ws.Range("I74").Value = 3160.2632
$ws.Range("K74").Value = 3160.2632
$ws.Range("M74").Value = -2286.2632
$ws.Range("H77").Value = 3764.0476
$ws.Range("I77").Value = 3160.2632
$ws.Range("K77").Value = 15801.316
$ws.Range("M77").Value = -11433.316
$ws.Range("H97").Value = 876.8570999999999
$ws.Range("I97").Value = 665.05554
$ws.Range("K97").Value = 665.05554
$ws.Range("M97").Value = -169.05554
$ws.Range("H132").Value = 4252.9536
$ws.Range("I132").Value = 3507.8857
$ws.Range("K132").Value = 10523.6571
$ws.Range("M132").Value = -7993.6571

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H9").Value = 24995
$ws.Range("J9").Value = 24995
$ws.Range("L9").Value = 24995
$ws.Range("N9").Value = -25331

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5195.7407
$ws.Range("I31").Value = 6572.857
$ws.Range("J31").Value = 4713.75
$ws.Range("K31").Value = 6572.857
$ws.Range("L31").Value = 4713.75
$ws.Range("M31").Value = -6277.857
$ws.Range("N31").Value = -5303.75
$ws.Range("H33").Value = 25995
$ws.Range("J33").Value = 25995
$ws.Range("L33").Value = 25995
$ws.Range("N33").Value = -26753
$ws.Range("H34").Value = 5195.7407
$ws.Range("I34").Value = 6572.857
$ws.Range("J34").Value = 4713.75
$ws.Range("K34").Value = 6572.857
$ws.Range("L34").Value = 4713.75
$ws.Range("M34").Value = -6370.857
$ws.Range("N34").Value = -5117.75
$ws.Range("H58").Value = 8679.866
$ws.Range("J58").Value = 9549.9
$ws.Range("L58").Value = 9549.9
$ws.Range("N58").Value = -9955.9
$ws.Range("H105").Value = 2292.0833
$ws.Range("I105").Value = 2292.0833
$ws.Range("K105").Value = 2292.0833
$ws.Range("M105").Value = -545.0832999999998
$ws.Range("H132").Value = 2828.0715
$ws.Range("I132").Value = 2828.0715
$ws.Range("K132").Value = 8484.2145
$ws.Range("M132").Value = -5954.2145
$ws.Range("H134").Value = 5569.643
$ws.Range("I134").Value = 4453.4707
$ws.Range("K134").Value = 13360.4121
$ws.Range("M134").Value = -10825.4121
$ws.Range("H136").Value = 8679.866
$ws.Range("J136").Value = 9549.9
$ws.Range("L136").Value = 28649.7
$ws.Range("N136").Value = -33749.7
$ws.Range("H140").Value = 117856
$ws.Range("J140").Value = 117856
$ws.Range("L140").Value = 117856
$ws.Range("N140").Value = -128216

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 11044.75
$ws.Range("I3").Value = 8978.154
$ws.Range("K3").Value = 26934.462
$ws.Range("M3").Value = -26822.462
$ws.Range("H32").Value = 632.8333
$ws.Range("I32").Value = 399.25
$ws.Range("J32").Value = 1100
$ws.Range("K32").Value = 1197.75
$ws.Range("L32").Value = 3300
$ws.Range("M32").Value = -914.75
$ws.Range("N32").Value = -3866
$ws.Range("H81").Value = 8554
$ws.Range("I81").Value = 610.6667
$ws.Range("J81").Value = 13320
$ws.Range("K81").Value = 1832.0001
$ws.Range("L81").Value = 39960
$ws.Range("M81").Value = -709.0001
$ws.Range("N81").Value = -42206
$ws.Range("H84").Value = 8554
$ws.Range("I84").Value = 610.6667
$ws.Range("J84").Value = 13320
$ws.Range("K84").Value = 5496.0003
$ws.Range("L84").Value = 119880
$ws.Range("M84").Value = 119.9997000000003
$ws.Range("N84").Value = -131112
$ws.Range("H132").Value = 38462944
$ws.Range("J132").Value = 1875.1666
$ws.Range("L132").Value = 16876.4994
$ws.Range("N132").Value = -21936.4994
$ws.Range("H133").Value = 2399.8572
$ws.Range("I133").Value = 2299.8333
$ws.Range("J133").Value = 3000
$ws.Range("K133").Value = 6899.499899999999
$ws.Range("L133").Value = 9000
$ws.Range("M133").Value = -1839.499899999999
$ws.Range("N133").Value = -19120

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H74").Value = 41998
$ws.Range("I74").Value = 30000
$ws.Range("J74").Value = 44997.5
$ws.Range("K74").Value = 30000
$ws.Range("L74").Value = 44997.5
$ws.Range("M74").Value = -29064
$ws.Range("N74").Value = -46869.5
$ws.Range("H77").Value = 41998
$ws.Range("I77").Value = 30000
$ws.Range("J77").Value = 44997.5
$ws.Range("K77").Value = 90000
$ws.Range("L77").Value = 134992.5
$ws.Range("M77").Value = -85320
$ws.Range("N77").Value = -144352.5
$ws.Range("H80").Value = 4620
$ws.Range("I80").Value = 7480
$ws.Range("K80").Value = 7480
$ws.Range("M80").Value = -6482
$ws.Range("H83").Value = 4620
$ws.Range("I83").Value = 7480
$ws.Range("K83").Value = 37400
$ws.Range("M83").Value = -32408
$ws.Range("H132").Value = 4382.385
$ws.Range("I132").Value = 3914.25
$ws.Range("K132").Value = 11742.75
$ws.Range("M132").Value = -9212.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H94").Value = 4000
$ws.Range("I94").Value = 4000
$ws.Range("K94").Value = 4000
$ws.Range("M94").Value = -3324
$ws.Range("H122").Value = 3063.32
$ws.Range("I122").Value = 3063.32
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 9189.960000000001
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -6739.960000000001
$ws.Range("N122").Value = ""
$ws.Range("H132").Value = 7238.354
$ws.Range("I132").Value = 7003.171
$ws.Range("J132").Value = 8615.857
$ws.Range("K132").Value = 21009.513
$ws.Range("L132").Value = 25847.571
$ws.Range("M132").Value = -18479.513
$ws.Range("N132").Value = -30907.571
$ws.Range("H136").Value = 4615.5713
$ws.Range("I136").Value = 3849.625
$ws.Range("K136").Value = 11548.875
$ws.Range("M136").Value = -8998.875

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H76").Value = 39996.332
$ws.Range("J76").Value = 39996.332
$ws.Range("L76").Value = 39996.332
$ws.Range("N76").Value = -40626.332
$ws.Range("H79").Value = 39996.332
$ws.Range("J79").Value = 39996.332
$ws.Range("L79").Value = 39996.332
$ws.Range("N79").Value = -42180.332
$ws.Range("H113").Value = 1100.3334
$ws.Range("I113").Value = 1300.4
$ws.Range("K113").Value = 3901.2
$ws.Range("M113").Value = -1731.2
$ws.Range("H126").Value = 1734.7894
$ws.Range("I126").Value = 1734.7894
$ws.Range("K126").Value = 5204.3682
$ws.Range("M126").Value = -2734.3682
$ws.Range("H132").Value = 2640.5098
$ws.Range("I132").Value = 2256.7446
$ws.Range("K132").Value = 6770.2338
$ws.Range("M132").Value = -4240.2338
$ws.Range("H136").Value = 6455.0454
$ws.Range("I136").Value = 5399.3335
$ws.Range("K136").Value = 16198.0005
$ws.Range("M136").Value = -13648.0005
$ws.Range("H137").Value = 80000
$ws.Range("J137").Value = 80000
$ws.Range("L137").Value = 80000
$ws.Range("N137").Value = -90200
